# Add diversification measure plotting
# Adds three new config rows (18-20, Config# 17-19) to Sheet1, reusing
# existing shared strings for the repeated labels and introducing a new
# "Add TE to MV" comment string, then moves the selection to the new
# last row (A20) to mirror the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 18 (Config 17) ---------------------------------------------------
$ws.Range("A18").Value = 17
$ws.Range("C18").Value = 0.06274
$ws.Range("D18").Value = 0.1
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = "Yes"
$ws.Range("G18").Value = "Yes"
$ws.Range("H18").Value = 0.025
$ws.Range("I18").Value = "Add TE to MV"

# --- Row 19 (Config 18) ---------------------------------------------------
$ws.Range("A19").Value = 18
$ws.Range("C19").Value = 0.06274
$ws.Range("D19").Value = 0.1
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = "None"
$ws.Range("G19").Value = "Yes"
$ws.Range("H19").Value = 0.25
$ws.Range("I19").Value = "Add TE to MV"

# --- Row 20 (Config 19) ---------------------------------------------------
$ws.Range("A20").Value = 19
$ws.Range("C20").Value = 0.06274
$ws.Range("D20").Value = 0.1
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = "None"
$ws.Range("G20").Value = "Yes"
$ws.Range("H20").Value = 0.25
$ws.Range("I20").Value = "Add TE to MV"

# Move / leave the selection on the newly added last row, same as the
# author's saved cursor position.
$ws.Range("A20").Select()
